# Applies the scheduled-runner price/profit refresh to the Zodiark_Profits workbook.
# Each worksheet (one per crafting job) gets updated currentAveragePrice /
# LevePrice / LeveProfit figures in columns H-N for the affected rows.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4191.0835
$ws.Range("I2").Value = 585.1111
$ws.Range("K2").Value = 585.1111
$ws.Range("M2").Value = -472.1111
$ws.Range("H41").Value = 240.94737
$ws.Range("I41").Value = 146.42857
$ws.Range("J41").Value = 296.08334
$ws.Range("K41").Value = 146.42857
$ws.Range("L41").Value = 296.08334
$ws.Range("M41").Value = 293.57143
$ws.Range("N41").Value = -1176.08334
$ws.Range("H43").Value = 7665
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 16995
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 16995
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -17133
$ws.Range("H86").Value = 1908.7
$ws.Range("J86").Value = 1431
$ws.Range("L86").Value = 1431
$ws.Range("N86").Value = -3677
$ws.Range("H89").Value = 1908.7
$ws.Range("J89").Value = 1431
$ws.Range("L89").Value = 7155
$ws.Range("N89").Value = -18387
$ws.Range("H107").Value = 460.34616
$ws.Range("I107").Value = 461.25
$ws.Range("K107").Value = 461.25
$ws.Range("M107").Value = 1458.75
$ws.Range("H114").Value = 67000
$ws.Range("J114").Value = 67000
$ws.Range("L114").Value = 67000
$ws.Range("N114").Value = -75678
$ws.Range("H116").Value = 5457.75
$ws.Range("J116").Value = 5462
$ws.Range("L116").Value = 5462
$ws.Range("N116").Value = -12346
$ws.Range("H138").Value = 1918.7637
$ws.Range("J138").Value = 2237.1843
$ws.Range("L138").Value = 6711.5529
$ws.Range("N138").Value = -16991.5529

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2873
$ws.Range("J74").Value = 3228.3215
$ws.Range("L74").Value = 3228.3215
$ws.Range("N74").Value = -4976.3215
$ws.Range("H76").Value = 22414
$ws.Range("J76").Value = 22414
$ws.Range("L76").Value = 22414
$ws.Range("N76").Value = -23090
$ws.Range("H77").Value = 2873
$ws.Range("J77").Value = 3228.3215
$ws.Range("L77").Value = 16141.6075
$ws.Range("N77").Value = -24877.6075
$ws.Range("H79").Value = 22414
$ws.Range("J79").Value = 22414
$ws.Range("L79").Value = 22414
$ws.Range("N79").Value = -24754
$ws.Range("H101").Value = 135258.4
$ws.Range("J101").Value = 135258.4
$ws.Range("L101").Value = 135258.4
$ws.Range("N101").Value = -141748.4
$ws.Range("H109").Value = 22500.4
$ws.Range("J109").Value = 22500.4
$ws.Range("L109").Value = 22500.4
$ws.Range("N109").Value = -25274.4
$ws.Range("H119").Value = 49750
$ws.Range("J119").Value = 49750
$ws.Range("L119").Value = 49750
$ws.Range("N119").Value = -59426
$ws.Range("H124").Value = 23539.25
$ws.Range("J124").Value = 23539.25
$ws.Range("L124").Value = 23539.25
$ws.Range("N124").Value = -33359.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 705.5
$ws.Range("I94").Value = 392
$ws.Range("K94").Value = 392
$ws.Range("M94").Value = 59
$ws.Range("H100").Value = 19697.334
$ws.Range("J100").Value = 19697.334
$ws.Range("L100").Value = 19697.334
$ws.Range("N100").Value = -21861.334
$ws.Range("H105").Value = 1612.8823
$ws.Range("I105").Value = 1588.6875
$ws.Range("K105").Value = 1588.6875
$ws.Range("M105").Value = 158.3125

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 14784.833
$ws.Range("J88").Value = 14784.833
$ws.Range("L88").Value = 14784.833
$ws.Range("N88").Value = -15596.833
$ws.Range("H91").Value = 14784.833
$ws.Range("J91").Value = 14784.833
$ws.Range("L91").Value = 14784.833
$ws.Range("N91").Value = -17592.833
$ws.Range("H109").Value = 32314.5
$ws.Range("J109").Value = 29666.334
$ws.Range("L109").Value = 29666.334
$ws.Range("N109").Value = -31746.334

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 62501050
$ws.Range("I4").Value = 125001096
$ws.Range("K4").Value = 375003288
$ws.Range("M4").Value = -375003176
$ws.Range("H87").Value = 35968.43
$ws.Range("I87").Value = 30556
$ws.Range("K87").Value = 91668
$ws.Range("M87").Value = -90420
$ws.Range("H90").Value = 35968.43
$ws.Range("I90").Value = 30556
$ws.Range("K90").Value = 275004
$ws.Range("M90").Value = -268764

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5081.6113
$ws.Range("I61").Value = 5233.5293
$ws.Range("J61").Value = 2499
$ws.Range("K61").Value = 5233.5293
$ws.Range("L61").Value = 2499
$ws.Range("M61").Value = -5031.5293
$ws.Range("N61").Value = -2903
$ws.Range("H93").Value = 3716.9614
$ws.Range("I93").Value = 3843.7058
$ws.Range("J93").Value = 3477.5557
$ws.Range("K93").Value = 3843.7058
$ws.Range("L93").Value = 3477.5557
$ws.Range("M93").Value = -2595.7058
$ws.Range("N93").Value = -5973.5557
$ws.Range("H104").Value = 16307
$ws.Range("J104").Value = 16307
$ws.Range("L104").Value = 16307
$ws.Range("N104").Value = -23295
$ws.Range("H113").Value = 5081.6113
$ws.Range("I113").Value = 5233.5293
$ws.Range("J113").Value = 2499
$ws.Range("K113").Value = 5233.5293
$ws.Range("L113").Value = 2499
$ws.Range("M113").Value = -3063.5293
$ws.Range("N113").Value = -6839
$ws.Range("H127").Value = 93639.60000000001
$ws.Range("J127").Value = 93639.60000000001
$ws.Range("L127").Value = 93639.60000000001
$ws.Range("N127").Value = -103559.6
$ws.Range("H132").Value = 4294.0654
$ws.Range("I132").Value = 4157.697
$ws.Range("K132").Value = 12473.091
$ws.Range("M132").Value = -9943.091
$ws.Range("H136").Value = 1848.76
$ws.Range("I136").Value = 1508.9166
$ws.Range("J136").Value = 10005
$ws.Range("K136").Value = 4526.7498
$ws.Range("L136").Value = 30015
$ws.Range("M136").Value = -1976.7498
$ws.Range("N136").Value = -35115
$ws.Range("H137").Value = 79999
$ws.Range("J137").Value = 79999
$ws.Range("L137").Value = 79999
$ws.Range("N137").Value = -90199
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360
$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 150000
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 100000
$ws.Range("J82").Value = 100000
$ws.Range("L82").Value = 100000
$ws.Range("N82").Value = -100766
$ws.Range("H85").Value = 100000
$ws.Range("J85").Value = 100000
$ws.Range("L85").Value = 100000
$ws.Range("N85").Value = -102652
$ws.Range("H103").Value = 35888.2
$ws.Range("J103").Value = 35888.2
$ws.Range("L103").Value = 35888.2
$ws.Range("N103").Value = -38232.2
$ws.Range("H107").Value = 239
$ws.Range("I107").Value = 212.625
$ws.Range("J107").Value = 379.66666
$ws.Range("K107").Value = 637.875
$ws.Range("L107").Value = 1138.99998
$ws.Range("M107").Value = 1282.125
$ws.Range("N107").Value = -4978.999980000001
$ws.Range("H122").Value = 7117.923
$ws.Range("I122").Value = 7140.091
$ws.Range("K122").Value = 21420.273
$ws.Range("M122").Value = -18970.273
$ws.Range("H132").Value = 1973.3572
$ws.Range("I132").Value = 1768.6666
$ws.Range("K132").Value = 5305.9998
$ws.Range("M132").Value = -2775.9998
$ws.Range("H136").Value = 2481.2415
$ws.Range("I136").Value = 2761
$ws.Range("K136").Value = 8283
$ws.Range("M136").Value = -5733
